# Update sheet name and title to reflect new "through" date (07-09 -> 07-10)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-07-10"
$ws.Range("I1").Value = "2022 (through 07-10)"

# Update the updated 2022 values in column I (through-date column)
$ws.Range("I5").Value = 114    # May
$ws.Range("I8").Value = 53     # August
$ws.Range("I14").Value = 858   # Total
